# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N) across the ALC, ARM, BSM,
# CRP, CUL, GSM and LTW leve sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 649.5
$ws.Range("I4").Value = 299
$ws.Range("K4").Value = 299
$ws.Range("M4").Value = -185

# Row 19
$ws.Range("H19").Value = 1971.3846
$ws.Range("I19").Value = 2245.111
$ws.Range("J19").Value = 1355.5
$ws.Range("K19").Value = 2245.111
$ws.Range("L19").Value = 1355.5
$ws.Range("M19").Value = -2070.111
$ws.Range("N19").Value = -1705.5

# Row 21
$ws.Range("H21").Value = 35021
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 35021
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 35021
$ws.Range("M21").Value = $null
$ws.Range("N21").Value = -35957

# Row 23
$ws.Range("H23").Value = 35021
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35021
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35021
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = -35489

# Row 29
$ws.Range("H29").Value = 3904.6
$ws.Range("J29").Value = 3906.75
$ws.Range("L29").Value = 11720.25
$ws.Range("N29").Value = -12282.25

# Row 43
$ws.Range("H43").Value = 3830.75
$ws.Range("I43").Value = 4441.6665
$ws.Range("K43").Value = 4441.6665
$ws.Range("M43").Value = -4372.6665

# Row 112
$ws.Range("H112").Value = 1748.3889
$ws.Range("J112").Value = 1860.6875
$ws.Range("L112").Value = 5582.0625
$ws.Range("N112").Value = -7798.0625

# Row 123
$ws.Range("H123").Value = 180000
$ws.Range("J123").Value = 180000
$ws.Range("L123").Value = 180000
$ws.Range("N123").Value = -189800

# Row 135
$ws.Range("H135").Value = 2393
$ws.Range("I135").Value = 2394.5
$ws.Range("K135").Value = 21550.5
$ws.Range("M135").Value = -19015.5

# Row 137
$ws.Range("H137").Value = 1622.1875
$ws.Range("I137").Value = 1186.591
$ws.Range("J137").Value = 2580.5
$ws.Range("K137").Value = 3559.773
$ws.Range("L137").Value = 7741.5
$ws.Range("M137").Value = -1009.773
$ws.Range("N137").Value = -12841.5

# Row 138
$ws.Range("H138").Value = 3264.5283
$ws.Range("I138").Value = 2676.4443
$ws.Range("J138").Value = 3384.818
$ws.Range("K138").Value = 8029.3329
$ws.Range("L138").Value = 10154.454
$ws.Range("M138").Value = -2889.3329
$ws.Range("N138").Value = -20434.454

# Row 141
$ws.Range("H141").Value = 5937.5
$ws.Range("I141").Value = 5937.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 17812.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -12632.5
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2002.2727
$ws.Range("I2").Value = 2059.5557
$ws.Range("K2").Value = 2059.5557
$ws.Range("M2").Value = -1946.5557

# Row 102
$ws.Range("H102").Value = 1381.909
$ws.Range("I102").Value = 1381.909
$ws.Range("K102").Value = 1381.909
$ws.Range("M102").Value = 240.0909999999999

# Row 116
$ws.Range("H116").Value = 2002.2727
$ws.Range("I116").Value = 2059.5557
$ws.Range("K116").Value = 2059.5557
$ws.Range("M116").Value = 234.4443000000001

# Row 122
$ws.Range("H122").Value = 8800
$ws.Range("I122").Value = 8800
$ws.Range("K122").Value = 26400
$ws.Range("M122").Value = -23950

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2002.2727
$ws.Range("I3").Value = 2059.5557
$ws.Range("K3").Value = 2059.5557
$ws.Range("M3").Value = -1945.5557

# Row 94
$ws.Range("H94").Value = 1327.4546
$ws.Range("J94").Value = 2586.25
$ws.Range("L94").Value = 2586.25
$ws.Range("N94").Value = -3488.25

# Row 105
$ws.Range("H105").Value = 4501.875
$ws.Range("I105").Value = 4073.5715
$ws.Range("K105").Value = 4073.5715
$ws.Range("M105").Value = -2326.5715

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2243.3333
$ws.Range("I16").Value = 2242.5
$ws.Range("K16").Value = 2242.5
$ws.Range("M16").Value = -1955.5

# Row 107
$ws.Range("H107").Value = 2595.4614
$ws.Range("I107").Value = 1750
$ws.Range("J107").Value = 2665.9167
$ws.Range("K107").Value = 1750
$ws.Range("L107").Value = 2665.9167
$ws.Range("M107").Value = 170
$ws.Range("N107").Value = -6505.9167

# Row 113
$ws.Range("H113").Value = 2243.3333
$ws.Range("I113").Value = 2242.5
$ws.Range("K113").Value = 2242.5
$ws.Range("M113").Value = -72.5

# Row 134
$ws.Range("H134").Value = 3456
$ws.Range("I134").Value = 2829.6
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 8488.799999999999
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -5953.799999999999
$ws.Range("N134").Value = -18570

# Row 141
$ws.Range("H141").Value = 64206.285
$ws.Range("I141").Value = 100000
$ws.Range("J141").Value = 58240.668
$ws.Range("K141").Value = 100000
$ws.Range("L141").Value = 58240.668
$ws.Range("M141").Value = -94820
$ws.Range("N141").Value = -68600.66800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 499
$ws.Range("I8").Value = 499
$ws.Range("K8").Value = 1497
$ws.Range("M8").Value = -1358

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null

# Row 40
$ws.Range("H40").Value = 199.5
$ws.Range("I40").Value = 199.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 798
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -729
$ws.Range("N40").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1726.875
$ws.Range("I102").Value = 1552.6666
$ws.Range("K102").Value = 1552.6666
$ws.Range("M102").Value = 69.33339999999998

# Row 113
$ws.Range("H113").Value = 1378.6666
$ws.Range("I113").Value = 1378.6666
$ws.Range("K113").Value = 1378.6666
$ws.Range("M113").Value = 791.3334

# Row 122
$ws.Range("H122").Value = 2078.25
$ws.Range("I122").Value = 1437.6666
$ws.Range("K122").Value = 4312.9998
$ws.Range("M122").Value = -1862.9998

# Row 126
$ws.Range("H126").Value = 1174.5
$ws.Range("I126").Value = 1199
$ws.Range("J126").Value = 1150
$ws.Range("K126").Value = 3597
$ws.Range("L126").Value = 3450
$ws.Range("M126").Value = -1127
$ws.Range("N126").Value = -8390

# Row 133
$ws.Range("H133").Value = 70236.336
$ws.Range("I133").Value = 60709
$ws.Range("K133").Value = 60709
$ws.Range("M133").Value = -55649

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null

# Row 82
$ws.Range("H82").Value = 1330.3334
$ws.Range("I82").Value = 995.5
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 995.5
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -634.5
$ws.Range("N82").Value = -2722

# Row 85
$ws.Range("H85").Value = 1330.3334
$ws.Range("I85").Value = 995.5
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 995.5
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = 252.5
$ws.Range("N85").Value = -4496

# Row 122
$ws.Range("H122").Value = 5500
$ws.Range("I122").Value = 5500
$ws.Range("K122").Value = 16500
$ws.Range("M122").Value = -14050

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
